$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TC01_Verify_HomePage")
$ws2 = $wb.Worksheets.Item("Testdata")

# Rename "HomeCarousel" -> "HeroBanner" everywhere it is used as test data.
$ws1.Range("C4").Value = "HeroBanner"
$ws1.Range("E4").Value = "HeroBanner"
$ws2.Range("A4").Value = "HeroBanner"

# Update sheet1 selection (no longer the active tab).
$ws1.Activate()
$ws1.Range("B3:E5").Select()

# Update sheet2 selection and make it the active tab (saved last/active).
$ws2.Activate()
$ws2.Range("B11").Select()
